$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for rows on both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, which carry
# identical event listings.
$updates = @{
    4  = 4681
    7  = 1406
    10 = 1188
    11 = 30
    12 = 661
    13 = 60
    14 = 51
    16 = 282
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
